# upload sensor board Robin + BOM
# Adds a "Resistors (sensorboard)" / "Capacitors (sensorboards)" / "Thermistor"
# component block to the BOM sheet (rows 15-24), removes the now-unused
# blank row 33, tweaks a couple of leftover placeholder rows, widens
# column A, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Set-RowFormats($row) {
    # D column: "0000" part-number style numeric format
    $ws.Range("D$row").NumberFormat = "0000"

    # F column: reuse the existing currency ("<euro> #,##0.00") style
    $ws.Range("F14").Copy()
    $ws.Range("F$row").PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Row 15 - Thermistor
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Thermistor"
$ws.Range("B15").Value = "100kOhm thermistor"
$ws.Range("C15").Value = "NCU18WF104D60RB"
$ws.Range("D15").Value = 603
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.297
$ws.Range("G15").Value = "https://www.mouser.be/ProductDetail/Murata-Electronics/NCU18WF104D60RB?qs=sGAEpiMZZMuBd0%252BwiCVS21gZfQ6Dyzsfx0RadtHN9DipnknzDvt5hw%3D%3D"
Set-RowFormats 15
$ws.Range("G10").Copy()
$ws.Range("G15").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 16 - Resistors (sensorboard) header row / 100 kOhm
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Resistors (sensorboard)"
$ws.Range("B16").Value = "100 kOhm"
$ws.Range("C16").Value = "CRCW0603100KFKEAC"
$ws.Range("C4").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("D16").Value = 603
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 0.09
$ws.Range("G16").Value = "https://www.mouser.be/ProductDetail/Vishay-Dale/CRCW0603100KFKEAC?qs=sGAEpiMZZMtlubZbdhIBIIZe04wfiaJWGPWKSQhf9Xo%3D"
Set-RowFormats 16
$ws.Range("G4").Copy()
$ws.Range("G16").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 17 - 10 kOhm
# ---------------------------------------------------------------------
$ws.Range("B17").Value = "10 kOhm"
$ws.Range("C17").Value = "CRCW060310K0FKEAC"
$ws.Range("C5").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("D17").Value = 603
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 0.135
$ws.Range("G17").Value = "https://www.mouser.be/ProductDetail/Vishay-Dale/CRCW060310K0FKEAC?qs=sGAEpiMZZMtlubZbdhIBIIZe04wfiaJWNE%252B7tlPkrYc%3D"
Set-RowFormats 17
$ws.Range("G10").Copy()
$ws.Range("G17").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 18 - 4,7 kOhm
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "4,7 kOhm"
$ws.Range("C18").Value = "CRCW06034K70FKEAC"
$ws.Range("D18").Value = 603
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.09
$ws.Range("G18").Value = "https://www.mouser.be/ProductDetail/Vishay-Dale/CRCW06034K70FKEAC?qs=sGAEpiMZZMtlubZbdhIBIIZe04wfiaJWtK1b03yAW%2Fw%3D"
Set-RowFormats 18
$ws.Range("G10").Copy()
$ws.Range("G18").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 19 - 22 Ohm
# ---------------------------------------------------------------------
$ws.Range("B19").Value = "22 Ohm"
$ws.Range("C19").Value = "CRCW060322R0FKEAC"
$ws.Range("D19").Value = 603
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.09
$ws.Range("G19").Value = "https://www.mouser.be/ProductDetail/Vishay-Dale/CRCW060322R0FKEAC?qs=sGAEpiMZZMtlubZbdhIBIIZe04wfiaJWhp35UV2eBC0%3D"
Set-RowFormats 19
$ws.Range("G10").Copy()
$ws.Range("G19").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 20 - Capacitors (sensorboards) header row / 100 nF
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Capacitors (sensorboards)"
$ws.Range("B20").Value = "100 nF"
$ws.Range("C20").Value = "GCJ188R71C104MA01D"
$ws.Range("D20").Value = 603
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.171
$ws.Range("G20").Value = "https://www.mouser.be/ProductDetail/Murata-Electronics/GCJ188R71C104MA01D?qs=sGAEpiMZZMs0AnBnWHyRQID2xuQsFd1GofHLF%2FZjZNHccNCeYysMBg%3D%3D"
Set-RowFormats 20
$ws.Range("G10").Copy()
$ws.Range("G20").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 21 - 1 uF
# ---------------------------------------------------------------------
$ws.Range("B21").Value = "1 µF"
$ws.Range("C21").Value = "GRM188R61C105KA12D"
$ws.Range("D21").Value = 603
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.144
$ws.Range("G21").Value = "https://www.mouser.be/ProductDetail/Murata-Electronics/GRM188R61C105KA12D?qs=sGAEpiMZZMs0AnBnWHyRQOK8EV8n4V50ILtJpoKun1k%3D"
Set-RowFormats 21
$ws.Range("G10").Copy()
$ws.Range("G21").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 22 - 2,2 uF
# ---------------------------------------------------------------------
$ws.Range("B22").Value = "2,2 µF"
$ws.Range("C22").Value = "GRM188C71E225KE11J"
$ws.Range("D22").Value = 603
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.279
$ws.Range("G22").Value = "https://www.mouser.be/ProductDetail/Murata-Electronics/GRM188C71E225KE11J?qs=sGAEpiMZZMs0AnBnWHyRQN7%2FAA2D2lPPoIBVQxy4%252BIYlGUALURHHzw%3D%3D"
Set-RowFormats 22
$ws.Range("G10").Copy()
$ws.Range("G22").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 23 - 4,7 uF
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "4,7 µF"
$ws.Range("C23").Value = "GRM188C71A475KE21D"
$ws.Range("D23").Value = 603
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.29
$ws.Range("G23").Value = "https://www.mouser.be/ProductDetail/Murata-Electronics/GRM188C71A475KE21D?qs=sGAEpiMZZMs0AnBnWHyRQCZFsEygxoaDUfTwU11CgLoHay4CQrnjdQ%3D%3D"
Set-RowFormats 23
$ws.Range("G10").Copy()
$ws.Range("G23").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 24 - 10 uF
# ---------------------------------------------------------------------
$ws.Range("B24").Value = "10 µF"
$ws.Range("C24").Value = "GRM188C80J106MEA6D"
$ws.Range("D24").Value = 603
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.171
$ws.Range("G24").Value = "https://www.mouser.be/ProductDetail/Murata-Electronics/GRM188C80J106MEA6D?qs=sGAEpiMZZMs0AnBnWHyRQN7%2FAA2D2lPPHKNFzVW0UbJeWppMSPj27Q%3D%3D"
Set-RowFormats 24
$ws.Range("G10").Copy()
$ws.Range("G24").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Leftover placeholder rows: 25/26 lose their currency format (now plain),
# and row 33 (trailing blank placeholder) is removed entirely.
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("F25").PasteSpecial($xlPasteFormats)
$ws.Range("F26").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(33).Delete()

# ---------------------------------------------------------------------
# Column A is widened to fit the new, longer part names.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24.1666666666667

# ---------------------------------------------------------------------
# Move the active selection, matching the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("B11").Select()
